$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '46.547.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +6.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.293.00'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +3.16%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '304.43'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.58%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '100.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +11.68%  '

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +2.49%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.73%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.67'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +11.03%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0790'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.88%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.35'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +5.93%  '

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.25%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.640.58'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +3.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.291.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.90'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +3.58%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +4.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '46.540.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +6.40%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.07'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +13.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0937'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.88%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.03'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +1.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '66.28'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '247.43'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.72%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.91'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.39%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.999'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.92'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '43.22'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +12.24%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.39%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.83'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +5.13%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.93'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +3.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.80'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +12.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.64'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.24%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '146.85'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -4.11%  '

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +5.06%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.23'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +13.92%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.116'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +12.68%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.04%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.60%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '15.97'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +21.09%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.06'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +12.16%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.38'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +5.94%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.40%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.98'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +11.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.832.30'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.03%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '87.18'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +19.63%  '

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +7.31%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '73.28'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +7.94%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.89'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +8.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '95.68'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.46%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.518.30'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +3.04%  '
